$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "236.49", "1.00").
# Excel auto-converts such strings to numeric values on assignment, which
# would silently drop significant trailing zeros / exact formatting (e.g.
# "1.00" -> 1, "5.50" -> 5.5). The source data stores these as text, so we
# format the cell as Text first, assign the literal value, then restore the
# Normal style so no stray number format lingers on the cell.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.49'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '617.13'
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.368'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.735'
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000244'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.82'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.50'
$ws.Range("D15").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.65'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.93'
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '439.51'
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.93'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '87.43'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.73'
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.123'
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.225'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.24'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.170'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.60'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '25.98'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '501.94'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.35'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.92'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.443'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.68'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.39'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.11'
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '158.28'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.705'
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '43.92'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.39'
$ws.Range("D51").Style = "Normal"

# Remaining cell updates (values that round-trip safely as plain text)
$ws.Range("D2").Value = '90.721.13'
$ws.Range("E2").Value = '  +0.64%  '
$ws.Range("D3").Value = '3.130.35'
$ws.Range("E3").Value = '  +0.88%  '
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("E6").Value = '  -0.66%  '
$ws.Range("E7").Value = '  +4.63%  '
$ws.Range("E8").Value = '  +1.82%  '
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("E10").Value = '  +2.36%  '
$ws.Range("D11").Value = '2.388.81'
$ws.Range("E11").Value = '  -22.95%  '
$ws.Range("E12").Value = '  +2.77%  '
$ws.Range("E13").Value = '  -2.53%  '
$ws.Range("E14").Value = '  -1.64%  '
$ws.Range("E15").Value = '  +2.01%  '
$ws.Range("D16").Value = '90.499.47'
$ws.Range("E16").Value = '  +0.69%  '
$ws.Range("D17").Value = '3.699.02'
$ws.Range("E17").Value = '  +0.81%  '
$ws.Range("D18").Value = '3.117.34'
$ws.Range("E18").Value = '  +0.72%  '
$ws.Range("E19").Value = '  -5.43%  '
$ws.Range("E20").Value = '  +7.54%  '
$ws.Range("E21").Value = '  +6.77%  '
$ws.Range("E22").Value = '  -6.18%  '
$ws.Range("E23").Value = '  +1.30%  '
$ws.Range("E24").Value = '  +1.73%  '
$ws.Range("E25").Value = '  +5.56%  '
$ws.Range("E26").Value = '  +1.07%  '
$ws.Range("E27").Value = '  -1.01%  '
$ws.Range("D28").Value = '3.296.00'
$ws.Range("E28").Value = '  +0.69%  '
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("E30").Value = '  +40.84%  '
$ws.Range("E31").Value = '  +6.77%  '
$ws.Range("E32").Value = '  +16.51%  '
$ws.Range("E33").Value = '  +1.09%  '
$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("E34").Value = '  +0.27%  '
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("E35").Value = '  +11.79%  '
$ws.Range("E36").Value = '  +5.78%  '
$ws.Range("E37").Value = '  +1.24%  '
$ws.Range("E38").Value = '  +0.73%  '
$ws.Range("B39").Value = 'Fetch.AI'
$ws.Range("C39").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("E39").Value = '  +6.31%  '
$ws.Range("B40").Value = 'PancakeSwap'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("E40").Value = '  +1.93%  '
$ws.Range("E41").Value = '  +11.29%  '
$ws.Range("E42").Value = '  +1.04%  '
$ws.Range("E43").Value = '  -9.80%  '
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("E46").Value = '  +4.71%  '
$ws.Range("E47").Value = '  +3.42%  '
$ws.Range("E48").Value = '  +1.72%  '
$ws.Range("E49").Value = '  +3.60%  '
$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("E50").Value = '  -1.22%  '
$ws.Range("B51").Value = 'Filecoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("E51").Value = '  +0.65%  '
